$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H101").Value = 1119.7142
$ws.Range("I101").Value = 416.83334
$ws.Range("J101").Value = 1646.875
$ws.Range("K101").Value = 1250.50002
$ws.Range("L101").Value = 4940.625
$ws.Range("M101").Value = 371.4999800000001
$ws.Range("N101").Value = -8184.625

$ws.Range("H137").Value = 3210.5
$ws.Range("I137").Value = 1640.7222
$ws.Range("J137").Value = 5565.1665
$ws.Range("K137").Value = 4922.1666
$ws.Range("L137").Value = 16695.4995
$ws.Range("M137").Value = -2372.1666
$ws.Range("N137").Value = -21795.4995

$ws.Range("H138").Value = 3423.3
$ws.Range("I138").Value = 1225
$ws.Range("J138").Value = 3751.7815
$ws.Range("K138").Value = 3675
$ws.Range("L138").Value = 11255.3445
$ws.Range("M138").Value = 1465
$ws.Range("N138").Value = -21535.3445

$ws.Range("H141").Value = 10092.458
$ws.Range("I141").Value = 11446
$ws.Range("J141").Value = 3324.75
$ws.Range("K141").Value = 34338
$ws.Range("L141").Value = 9974.25
$ws.Range("M141").Value = -29158
$ws.Range("N141").Value = -20334.25

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4756.1562
$ws.Range("I32").Value = 3670.8474
$ws.Range("J32").Value = 17562.8
$ws.Range("K32").Value = 3670.8474
$ws.Range("L32").Value = 17562.8
$ws.Range("M32").Value = -3383.8474
$ws.Range("N32").Value = -18136.8

$ws.Range("H61").Value = 1068.8788
$ws.Range("I61").Value = 822
$ws.Range("J61").Value = 1727.2222
$ws.Range("K61").Value = 822
$ws.Range("L61").Value = 1727.2222
$ws.Range("M61").Value = -610
$ws.Range("N61").Value = -2151.2222

$ws.Range("H63").Value = 11545183
$ws.Range("I63").Value = 12594200
$ws.Range("K63").Value = 12594200
$ws.Range("M63").Value = -12593514

$ws.Range("H66").Value = 11545183
$ws.Range("I66").Value = 12594200
$ws.Range("K66").Value = 62971000
$ws.Range("M66").Value = -62967568

$ws.Range("H74").Value = 2406.923
$ws.Range("I74").Value = 2422.476
$ws.Range("J74").Value = 2341.6
$ws.Range("K74").Value = 2422.476
$ws.Range("L74").Value = 2341.6
$ws.Range("M74").Value = -1548.476
$ws.Range("N74").Value = -4089.6

$ws.Range("H77").Value = 2406.923
$ws.Range("I77").Value = 2422.476
$ws.Range("J77").Value = 2341.6
$ws.Range("K77").Value = 12112.38
$ws.Range("L77").Value = 11708
$ws.Range("M77").Value = -7744.380000000001
$ws.Range("N77").Value = -20444

$ws.Range("H136").Value = 1068.8788
$ws.Range("I136").Value = 822
$ws.Range("J136").Value = 1727.2222
$ws.Range("K136").Value = 2466
$ws.Range("L136").Value = 5181.6666
$ws.Range("M136").Value = 84
$ws.Range("N136").Value = -10281.6666

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H35").Value = 28764.8
$ws.Range("J35").Value = 28764.8
$ws.Range("L35").Value = 28764.8
$ws.Range("N35").Value = -29384.8

$ws.Range("H134").Value = 2478.8333
$ws.Range("I134").Value = 1269.3235
$ws.Range("K134").Value = 3807.9705
$ws.Range("M134").Value = -1272.9705

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 12823369
$ws.Range("I31").Value = 1484.0454
$ws.Range("J31").Value = 29416396
$ws.Range("K31").Value = 1484.0454
$ws.Range("L31").Value = 29416396
$ws.Range("M31").Value = -1189.0454
$ws.Range("N31").Value = -29416986

$ws.Range("H34").Value = 12823369
$ws.Range("I34").Value = 1484.0454
$ws.Range("J34").Value = 29416396
$ws.Range("K34").Value = 1484.0454
$ws.Range("L34").Value = 29416396
$ws.Range("M34").Value = -1282.0454
$ws.Range("N34").Value = -29416800

$ws.Range("H58").Value = 1783.4744
$ws.Range("I58").Value = 1555.6522
$ws.Range("J58").Value = 3530.111
$ws.Range("K58").Value = 1555.6522
$ws.Range("L58").Value = 3530.111
$ws.Range("M58").Value = -1352.6522
$ws.Range("N58").Value = -3936.111

$ws.Range("H68").Value = 56708.125
$ws.Range("J68").Value = 56708.125
$ws.Range("L68").Value = 56708.125
$ws.Range("N68").Value = -58206.125

$ws.Range("H71").Value = 56708.125
$ws.Range("J71").Value = 56708.125
$ws.Range("L71").Value = 170124.375
$ws.Range("N71").Value = -177612.375

$ws.Range("H94").Value = 1523.2222
$ws.Range("J94").Value = 1693.3572
$ws.Range("L94").Value = 1693.3572
$ws.Range("N94").Value = -2595.3572

$ws.Range("H132").Value = 2947.1765
$ws.Range("I132").Value = 2535.3928
$ws.Range("J132").Value = 4868.8335
$ws.Range("K132").Value = 7606.178400000001
$ws.Range("L132").Value = 14606.5005
$ws.Range("M132").Value = -5076.178400000001
$ws.Range("N132").Value = -19666.5005

$ws.Range("H134").Value = 4178.95
$ws.Range("I134").Value = 4621.077
$ws.Range("J134").Value = 3357.8572
$ws.Range("K134").Value = 13863.231
$ws.Range("L134").Value = 10073.5716
$ws.Range("M134").Value = -11328.231
$ws.Range("N134").Value = -15143.5716

$ws.Range("H136").Value = 1783.4744
$ws.Range("I136").Value = 1555.6522
$ws.Range("J136").Value = 3530.111
$ws.Range("K136").Value = 4666.9566
$ws.Range("L136").Value = 10590.333
$ws.Range("M136").Value = -2116.9566
$ws.Range("N136").Value = -15690.333

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 952.3182
$ws.Range("I113").Value = 739.2353000000001
$ws.Range("J113").Value = 1676.8
$ws.Range("K113").Value = 2217.7059
$ws.Range("L113").Value = 5030.4
$ws.Range("M113").Value = -47.70589999999993
$ws.Range("N113").Value = -9370.4

$ws.Range("H131").Value = 6173675.5
$ws.Range("I131").Value = 71429120
$ws.Range("J131").Value = 863.0135
$ws.Range("K131").Value = 214287360
$ws.Range("L131").Value = 2589.0405
$ws.Range("M131").Value = -214282320
$ws.Range("N131").Value = -12669.0405

$ws.Range("H136").Value = 2971.1304
$ws.Range("I136").Value = 2735.5293
$ws.Range("J136").Value = 3638.6667
$ws.Range("K136").Value = 8206.5879
$ws.Range("L136").Value = 10916.0001
$ws.Range("M136").Value = -3106.5879
$ws.Range("N136").Value = -21116.0001

$ws.Range("H138").Value = 2714.5715
$ws.Range("I138").Value = 2500.3333
$ws.Range("K138").Value = 7500.999899999999
$ws.Range("M138").Value = -2360.999899999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 4699.6
$ws.Range("I122").Value = 2570.8572
$ws.Range("J122").Value = 9666.666999999999
$ws.Range("K122").Value = 7712.571599999999
$ws.Range("L122").Value = 29000.001
$ws.Range("M122").Value = -5262.571599999999
$ws.Range("N122").Value = -33900.001

$ws.Range("H132").Value = 2336.182
$ws.Range("I132").Value = 1407.9259
$ws.Range("J132").Value = 6513.3335
$ws.Range("K132").Value = 4223.7777
$ws.Range("L132").Value = 19540.0005
$ws.Range("M132").Value = -1693.7777
$ws.Range("N132").Value = -24600.0005

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H81").Value = 69874.5
$ws.Range("J81").Value = 69874.5
$ws.Range("L81").Value = 69874.5
$ws.Range("N81").Value = -71870.5

$ws.Range("H82").Value = 1217.6444
$ws.Range("I82").Value = 649.6667
$ws.Range("J82").Value = 1866.762
$ws.Range("K82").Value = 649.6667
$ws.Range("L82").Value = 1866.762
$ws.Range("M82").Value = -288.6667
$ws.Range("N82").Value = -2588.762

$ws.Range("H84").Value = 69874.5
$ws.Range("J84").Value = 69874.5
$ws.Range("L84").Value = 209623.5
$ws.Range("N84").Value = -219607.5

$ws.Range("H85").Value = 1217.6444
$ws.Range("I85").Value = 649.6667
$ws.Range("J85").Value = 1866.762
$ws.Range("K85").Value = 649.6667
$ws.Range("L85").Value = 1866.762
$ws.Range("M85").Value = 598.3333
$ws.Range("N85").Value = -4362.762

$ws.Range("H93").Value = 4445934.5
$ws.Range("I93").Value = 9260326
$ws.Range("J93").Value = 1881.6154
$ws.Range("K93").Value = 9260326
$ws.Range("L93").Value = 1881.6154
$ws.Range("M93").Value = -9259078
$ws.Range("N93").Value = -4377.6154

$ws.Range("H132").Value = 11798.857
$ws.Range("I132").Value = 14531.579
$ws.Range("J132").Value = 8553.75
$ws.Range("K132").Value = 43594.737
$ws.Range("L132").Value = 25661.25
$ws.Range("M132").Value = -41064.737
$ws.Range("N132").Value = -30721.25

$ws.Range("H136").Value = 2386.205
$ws.Range("I136").Value = 1289.2812
$ws.Range("K136").Value = 3867.8436
$ws.Range("M136").Value = -1317.8436

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 11496936
$ws.Range("I132").Value = 1935.2941
$ws.Range("J132").Value = 27781520
$ws.Range("K132").Value = 5805.8823
$ws.Range("L132").Value = 83344560
$ws.Range("M132").Value = -3275.8823
$ws.Range("N132").Value = -83349620

$ws.Range("H136").Value = 1732.1833
$ws.Range("I136").Value = 565.53656
$ws.Range("J136").Value = 4249.684
$ws.Range("K136").Value = 1696.60968
$ws.Range("L136").Value = 12749.052
$ws.Range("M136").Value = 853.39032
$ws.Range("N136").Value = -17849.052
